# Auto-generated script to apply scheduled price/profit updates across Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H5").Value = 124.5
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 99
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 99
$ws.Range("M5").Value = -35
$ws.Range("N5").Value = -329
$ws.Range("H32").Value = 9497.857
$ws.Range("J32").Value = 10622.5
$ws.Range("L32").Value = 10622.5
$ws.Range("N32").Value = -11274.5
$ws.Range("H38").Value = 1210
$ws.Range("J38").Value = 6509
$ws.Range("L38").Value = 19527
$ws.Range("N38").Value = -20271
$ws.Range("H40").Value = 2466.6667
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2850
$ws.Range("H42").Value = 429.42856
$ws.Range("J42").Value = 403.66666
$ws.Range("L42").Value = 1210.99998
$ws.Range("N42").Value = -1670.99998
$ws.Range("H113").Value = 10000000
$ws.Range("I113").Value = 10000000
$ws.Range("K113").Value = 10000000
$ws.Range("M113").Value = -9996746
$ws.Range("H115").Value = 8359.799999999999
$ws.Range("I115").Value = 8359.799999999999
$ws.Range("K115").Value = 25079.4
$ws.Range("M115").Value = -23512.4
$ws.Range("H138").Value = 2654.762
$ws.Range("I138").Value = 1897.8334
$ws.Range("J138").Value = 2957.5334
$ws.Range("K138").Value = 5693.5002
$ws.Range("L138").Value = 8872.600199999999
$ws.Range("M138").Value = -553.5002000000004
$ws.Range("N138").Value = -19152.6002

$ws = $wb.Sheets("ARM")
$ws.Range("H2").Value = 878.3889
$ws.Range("I2").Value = 1082.0769
$ws.Range("J2").Value = 348.8
$ws.Range("K2").Value = 1082.0769
$ws.Range("L2").Value = 348.8
$ws.Range("M2").Value = -969.0769
$ws.Range("N2").Value = -574.8
$ws.Range("H45").Value = 2130.4285
$ws.Range("I45").Value = 2130.4285
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2130.4285
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1753.4285
$ws.Range("N45").ClearContents()
$ws.Range("H116").Value = 878.3889
$ws.Range("I116").Value = 1082.0769
$ws.Range("J116").Value = 348.8
$ws.Range("K116").Value = 1082.0769
$ws.Range("L116").Value = 348.8
$ws.Range("M116").Value = 1211.9231
$ws.Range("N116").Value = -4936.8

$ws = $wb.Sheets("BSM")
$ws.Range("H3").Value = 878.3889
$ws.Range("I3").Value = 1082.0769
$ws.Range("J3").Value = 348.8
$ws.Range("K3").Value = 1082.0769
$ws.Range("L3").Value = 348.8
$ws.Range("M3").Value = -968.0769
$ws.Range("N3").Value = -576.8
$ws.Range("H105").Value = 4034.3333
$ws.Range("I105").Value = 3601.125
$ws.Range("J105").Value = 7500
$ws.Range("K105").Value = 3601.125
$ws.Range("L105").Value = 7500
$ws.Range("M105").Value = -1854.125
$ws.Range("N105").Value = -10994
$ws.Range("H134").Value = 1660.7812
$ws.Range("I134").Value = 1551.8214
$ws.Range("K134").Value = 4655.4642
$ws.Range("M134").Value = -2120.4642

$ws = $wb.Sheets("CRP")
$ws.Range("H16").Value = 1456.6
$ws.Range("I16").Value = 1456.6
$ws.Range("K16").Value = 1456.6
$ws.Range("M16").Value = -1169.6
$ws.Range("H58").Value = 2802.3333
$ws.Range("I58").Value = 2291.5
$ws.Range("K58").Value = 2291.5
$ws.Range("M58").Value = -2088.5
$ws.Range("H105").Value = 1676.2222
$ws.Range("I105").Value = 1723.25
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1723.25
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 23.75
$ws.Range("N105").Value = -4794
$ws.Range("H107").Value = 1457.75
$ws.Range("I107").Value = 692.4
$ws.Range("J107").Value = 2733.3333
$ws.Range("K107").Value = 692.4
$ws.Range("L107").Value = 2733.3333
$ws.Range("M107").Value = 1227.6
$ws.Range("N107").Value = -6573.3333
$ws.Range("H113").Value = 1456.6
$ws.Range("I113").Value = 1456.6
$ws.Range("K113").Value = 1456.6
$ws.Range("M113").Value = 713.4000000000001
$ws.Range("H136").Value = 2802.3333
$ws.Range("I136").Value = 2291.5
$ws.Range("K136").Value = 6874.5
$ws.Range("M136").Value = -4324.5

$ws = $wb.Sheets("CUL")
$ws.Range("H5").Value = 1011.5
$ws.Range("I5").Value = 682
$ws.Range("K5").Value = 2046
$ws.Range("M5").Value = -1934
$ws.Range("H14").Value = 145
$ws.Range("I14").Value = 145
$ws.Range("K14").Value = 435
$ws.Range("M14").Value = -262
$ws.Range("H135").Value = 1011.5
$ws.Range("I135").Value = 682
$ws.Range("K135").Value = 6138
$ws.Range("M135").Value = -3603
$ws.Range("H137").Value = 3128.5715
$ws.Range("I137").Value = 2975
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 8925
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -3825
$ws.Range("N137").Value = -20199.9999

$ws = $wb.Sheets("GSM")
$ws.Range("H11").Value = 15005000
$ws.Range("I11").Value = 30000000
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 30000000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -29999861
$ws.Range("N11").Value = -10278
$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336
$ws.Range("H102").Value = 1136.7646
$ws.Range("J102").Value = 1191.5
$ws.Range("L102").Value = 1191.5
$ws.Range("N102").Value = -4435.5
$ws.Range("H132").Value = 3045.9375
$ws.Range("I132").Value = 1860
$ws.Range("J132").Value = 4570.7144
$ws.Range("K132").Value = 5580
$ws.Range("L132").Value = 13712.1432
$ws.Range("M132").Value = -3050
$ws.Range("N132").Value = -18772.1432

$ws = $wb.Sheets("LTW")
$ws.Range("H55").Value = 199
$ws.Range("I55").Value = 197.5
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 197.5
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = -24.5
$ws.Range("N55").Value = -546

$ws = $wb.Sheets("WVR")
$ws.Range("H100").Value = 1062.8
$ws.Range("I100").Value = 1140.5
$ws.Range("J100").Value = 752
$ws.Range("K100").Value = 2281
$ws.Range("L100").Value = 1504
$ws.Range("M100").Value = -1740
$ws.Range("N100").Value = -2586
$ws.Range("H113").Value = 1259.6154
$ws.Range("I113").Value = 1307.8182
$ws.Range("J113").Value = 994.5
$ws.Range("K113").Value = 3923.4546
$ws.Range("L113").Value = 2983.5
$ws.Range("M113").Value = -1753.4546
$ws.Range("N113").Value = -7323.5
$ws.Range("H122").Value = 3025.2307
$ws.Range("I122").Value = 2802.6365
$ws.Range("K122").Value = 8407.9095
$ws.Range("M122").Value = -5957.9095
$ws.Range("H136").Value = 2162.9666
$ws.Range("I136").Value = 1793.762
$ws.Range("J136").Value = 3024.4443
$ws.Range("K136").Value = 5381.286
$ws.Range("L136").Value = 9073.332900000001
$ws.Range("M136").Value = -2831.286
$ws.Range("N136").Value = -14173.3329
